# Update cryptocurrency price/volume data on the active worksheet.
# All D (Price) and E (Volume(1h)) columns hold text values (inline strings),
# so we must assign strings explicitly to avoid Excel reinterpreting them as
# numbers/percentages and altering their literal formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($cellRef, $value) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
}

# Row 2
Set-Text "D2" "278.85"
Set-Text "E2" "6.89%"

# Row 3
Set-Text "D3" "27.25"
Set-Text "E3" "0.32%"

# Row 4
Set-Text "D4" "4.797"
Set-Text "E4" "1.77%"

# Row 5
Set-Text "D5" "0.06246"
Set-Text "E5" "0.41%"

# Row 6
Set-Text "D6" "6.845"
Set-Text "E6" "1.69%"

# Row 7
Set-Text "D7" "0.8778"
Set-Text "E7" "3.05%"

# Row 8
Set-Text "D8" "0.9390"
Set-Text "E8" "3.22%"

# Row 9
Set-Text "D9" "0.1453"
Set-Text "E9" "3.71%"

# Row 10
Set-Text "D10" "0.05043"
Set-Text "E10" "4.63%"

# Row 11
Set-Text "D11" "0.07278"
Set-Text "E11" "2.79%"

# Row 12
Set-Text "D12" "0.03150"
Set-Text "E12" "0.75%"

# Row 13 (only D changes)
Set-Text "D13" "0.09035"

# Row 14
Set-Text "D14" "0.001548"
Set-Text "E14" "1.19%"

# Row 15
Set-Text "D15" "0.0006265"
Set-Text "E15" "1.79%"

# Row 16
Set-Text "D16" "0.006109"
Set-Text "E16" "0.69%"

# Row 17
Set-Text "D17" "3.468"
Set-Text "E17" "0.46%"

# Row 18
Set-Text "D18" "3.270"
Set-Text "E18" "3.12%"

# Row 19
Set-Text "D19" "2.255"
Set-Text "E19" "4.16%"

# Row 21
Set-Text "D21" "0.1310"
Set-Text "E21" "0.04%"

# Row 22
Set-Text "D22" "3.851"
Set-Text "E22" "-6.01%"

# Row 23
Set-Text "D23" "0.04331"
Set-Text "E23" "2.16%"

# Row 24
Set-Text "D24" "0.001175"
Set-Text "E24" "-3.73%"

# Row 25
Set-Text "D25" "0.004267"
Set-Text "E25" "4.52%"

# Row 26 (only E changes)
Set-Text "E26" "-0.16%"

# Row 27
Set-Text "D27" "0.0001613"
Set-Text "E27" "-1.60%"

# Row 40
Set-Text "D40" "0.04027"
Set-Text "E40" "2.94%"

# Row 41
Set-Text "D41" "0.006700"
Set-Text "E41" "62.98%"

# Row 42
Set-Text "D42" "0.1150"
Set-Text "E42" "3.56%"

# Row 43 (only E changes)
Set-Text "E43" "-1.11%"

# Row 44
Set-Text "D44" "0.01234"
Set-Text "E44" "-5.19%"

# Row 45
Set-Text "D45" "0.00005121"
Set-Text "E45" "-0.11%"

# Row 46 (only E changes)
Set-Text "E46" "-0.15%"

# Row 47 (only E changes)
Set-Text "E47" "2,958.24%"

# Row 48 (only E changes)
Set-Text "E48" "-12.18%"

# Row 49
Set-Text "D49" "0.00002098"
Set-Text "E49" "-0.15%"

# Row 50
Set-Text "D50" "0.0001998"
Set-Text "E50" "-0.15%"
